$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# duration option: type now also accepts a string, description documents the default
$ws.Range("B13").Value = "Number (milliseconds) \| String"

# rotateStyle option: type notation now escapes the pipe characters
$ws.Range("B10").Value = "Number (1\|2\|3)"

# scale option: type changed from ratio to float, description clarified
$ws.Range("B11").Value = "Number (float)"
$ws.Range("D11").Value = "Hovering over hoverTarget scales to this value. 1 = 100%"

# Update the CodePen example links for updateRate / translate / translateReverse / tilt / tiltReverse
$ws.Range("C3").Value = "[CodePen](https://codepen.io/maiCoding/pen/MZEMqv)"
$ws.Range("C4").Value = "[CodePen](https://codepen.io/maiCoding/pen/ZVaLOp)"
$ws.Range("C5").Value = "[CodePen](https://codepen.io/maiCoding/pen/ZVaLOp)"
$ws.Range("C6").Value = "[CodePen](https://codepen.io/maiCoding/pen/aPVQmw)"
$ws.Range("C7").Value = "[CodePen](https://codepen.io/maiCoding/pen/aPVQmw)"

# duration description now documents the default
$ws.Range("D13").Value = "How many milliseconds/seconds it takes for a transform transition to complete. Default is 200ms"

# update the active selection left on the sheet
$ws.Range("L15").Select()
